$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("O1").Value = "CONTA ORIGEM"
$ws.Range("P1").Value = "DESC CONTA ORIGEM"
$ws.Range("Q1").Value = "AQUISITION ORIGEM"

# New data cells (row 2)
$ws.Range("O2").Value = 17
$ws.Range("P2").Value = "INFO"

# Date value (2012-12-12 -> serial 41255), formatted as a short date (numFmtId 14)
$ws.Range("Q2").Value = 41255
$ws.Range("Q2").NumberFormat = "mm-dd-yy"

# Selection mirrors the new block that was just filled in
$ws.Range("O1:Q2").Select()
